$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5800
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 5500
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 5500
$ws.Range("M64").Value = -5752
$ws.Range("N64").Value = -5996

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5800
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 5500
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 5500
$ws.Range("M67").Value = -5142
$ws.Range("N67").Value = -7216

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1943.8889
$ws.Range("I98").Value = 2005.2941
$ws.Range("J98").Value = 900
$ws.Range("K98").Value = 2005.2941
$ws.Range("L98").Value = 900
$ws.Range("M98").Value = -507.2941000000001
$ws.Range("N98").Value = -3896

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2624
$ws.Range("I112").Value = 840
$ws.Range("J112").Value = 3584.6155
$ws.Range("K112").Value = 2520
$ws.Range("L112").Value = 10753.8465
$ws.Range("M112").Value = -1412
$ws.Range("N112").Value = -12969.8465

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1943.8889
$ws.Range("I122").Value = 2005.2941
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 6015.8823
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -3565.8823
$ws.Range("N122").Value = -7600

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1022.87756
$ws.Range("J129").Value = 1167.5
$ws.Range("L129").Value = 3502.5
$ws.Range("N129").Value = -13502.5

# ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 71833.336
$ws.Range("J139").Value = 71833.336
$ws.Range("L139").Value = 71833.336
$ws.Range("N139").Value = -82113.336

# ARM row 43
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 7980
$ws.Range("J43").Value = 7980
$ws.Range("L43").Value = 7980
$ws.Range("N43").Value = -8606

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1772.9387
$ws.Range("I61").Value = 1726.0286
$ws.Range("J61").Value = 1890.2142
$ws.Range("K61").Value = 1726.0286
$ws.Range("L61").Value = 1890.2142
$ws.Range("M61").Value = -1514.0286
$ws.Range("N61").Value = -2314.2142

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1039.2094
$ws.Range("I74").Value = 1061.3334
$ws.Range("J74").Value = 988.1539
$ws.Range("K74").Value = 1061.3334
$ws.Range("L74").Value = 988.1539
$ws.Range("M74").Value = -187.3334
$ws.Range("N74").Value = -2736.1539

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1039.2094
$ws.Range("I77").Value = 1061.3334
$ws.Range("J77").Value = 988.1539
$ws.Range("K77").Value = 5306.666999999999
$ws.Range("L77").Value = 4940.7695
$ws.Range("M77").Value = -938.6669999999995
$ws.Range("N77").Value = -13676.7695

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3085.8823
$ws.Range("I102").Value = 2964.2856
$ws.Range("J102").Value = 3653.3333
$ws.Range("K102").Value = 2964.2856
$ws.Range("L102").Value = 3653.3333
$ws.Range("M102").Value = -1342.2856
$ws.Range("N102").Value = -6897.3333

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1772.9387
$ws.Range("I136").Value = 1726.0286
$ws.Range("J136").Value = 1890.2142
$ws.Range("K136").Value = 5178.085800000001
$ws.Range("L136").Value = 5670.642599999999
$ws.Range("M136").Value = -2628.085800000001
$ws.Range("N136").Value = -10770.6426

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5914.4165
$ws.Range("I62").Value = 5664.778
$ws.Range("J62").Value = 6663.3335
$ws.Range("K62").Value = 5664.778
$ws.Range("L62").Value = 6663.3335
$ws.Range("M62").Value = -5040.778
$ws.Range("N62").Value = -7911.3335

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5914.4165
$ws.Range("I65").Value = 5664.778
$ws.Range("J65").Value = 6663.3335
$ws.Range("K65").Value = 28323.89
$ws.Range("L65").Value = 33316.6675
$ws.Range("M65").Value = -25203.89
$ws.Range("N65").Value = -39556.6675

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1105.2413
$ws.Range("I132").Value = 816.8148
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2450.4444
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = 79.55560000000014
$ws.Range("N132").Value = -20057

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 21740086
$ws.Range("I134").Value = 998.1818
$ws.Range("J134").Value = 500000000
$ws.Range("K134").Value = 2994.5454
$ws.Range("L134").Value = 1500000000
$ws.Range("M134").Value = -459.5454
$ws.Range("N134").Value = -1500005070

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2673956.8
$ws.Range("I2").Value = 4902160
$ws.Range("J2").Value = 112.9
$ws.Range("K2").Value = 29412960
$ws.Range("L2").Value = 677.4000000000001
$ws.Range("M2").Value = -29412847
$ws.Range("N2").Value = -903.4000000000001

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 44716.668
$ws.Range("I11").Value = 57221.43
$ws.Range("J11").Value = 950
$ws.Range("K11").Value = 171664.29
$ws.Range("L11").Value = 2850
$ws.Range("M11").Value = -171524.29
$ws.Range("N11").Value = -3130

# CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 259.33334
$ws.Range("J26").Value = 345
$ws.Range("L26").Value = 1035
$ws.Range("N26").Value = -1611

# CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 199.5
$ws.Range("I40").Value = 186.66667
$ws.Range("J40").Value = 212.33333
$ws.Range("K40").Value = 746.66668
$ws.Range("L40").Value = 849.33332
$ws.Range("M40").Value = -677.66668
$ws.Range("N40").Value = -987.33332

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1670
$ws.Range("I122").Value = 1483.3334
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 4450.0002
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -2000.0002
$ws.Range("N122").Value = -10750

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3687.4285
$ws.Range("I126").Value = 4162.4
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 12487.2
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -10017.2
$ws.Range("N126").Value = -12440

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3707.805
$ws.Range("I132").Value = 3660.276
$ws.Range("J132").Value = 3822.6667
$ws.Range("K132").Value = 10980.828
$ws.Range("L132").Value = 11468.0001
$ws.Range("M132").Value = -8450.828
$ws.Range("N132").Value = -16528.0001

# GSM row 141
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 59680
$ws.Range("J141").Value = 59680
$ws.Range("L141").Value = 59680
$ws.Range("N141").Value = -70040

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2496.8
$ws.Range("I122").Value = 2496
$ws.Range("K122").Value = 7488
$ws.Range("M122").Value = -5038

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2700.4888
$ws.Range("I132").Value = 1460.1818
$ws.Range("J132").Value = 3886.8696
$ws.Range("K132").Value = 4380.5454
$ws.Range("L132").Value = 11660.6088
$ws.Range("M132").Value = -1850.5454
$ws.Range("N132").Value = -16720.6088

# WVR row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 5793.875
$ws.Range("J45").Value = 5950.143
$ws.Range("L45").Value = 5950.143
$ws.Range("N45").Value = -6932.143

# WVR row 74
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8185.2
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 8981.5
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 8981.5
$ws.Range("M74").Value = -4064
$ws.Range("N74").Value = -10853.5

# WVR row 77
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 8185.2
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 8981.5
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 26944.5
$ws.Range("M77").Value = -10320
$ws.Range("N77").Value = -36304.5

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 938.88
$ws.Range("I136").Value = 812.9524
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 2438.8572
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = 111.1428000000001
$ws.Range("N136").Value = -9900
